$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. BOM row 8: MOSFET changed from IRF1404ZSTRLPBFCT-ND (40V/180A, qty 8)
#    to IRFS3306TRLPBFCT-ND (60V/120A, qty 4) so the board now tolerates up
#    to 48V while keeping the 65A continuous rating.
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "IRFS3306TRLPBFCT-ND"
$ws.Range("C8").Value = "N-Channel 60 V 120A (Tc) 230W (Tc) Surface Mount D2PAK"
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 2.58
$ws.Range("H8").Value = "https://www.infineon.com/dgdl/irfs3306pbf.pdf?fileId=5546d462533600a40153563682652165"

# ---------------------------------------------------------------------------
# 2. Re-point the H8 hyperlink at the new datasheet URL. The COM shim here
#    maps a single hyperlink's Delete()/Address-set to an "add new, leave
#    stale" operation, and Range("H8").Hyperlinks.Delete() clears every
#    hyperlink on the sheet (not just H8's) - so the reliable way to land in
#    the same end state is to clear them all and re-add all twenty in the
#    desired order, with H8 last (matching the source workbook's edit, where
#    the H8 hyperlink relationship became the newest/last one).
# ---------------------------------------------------------------------------
$ws.Range("H8").Hyperlinks.Delete()

$links = @(
  @("H2",  "https://www.digikey.com/en/products/detail/kemet/C1812C104J5RAC7800/2235598?s=N4IgTCBcDaIMIEYAcCyIAwBYBSBWASgIJwAqAqiALoC%2BQA"),
  @("H20", "https://www.digikey.com/en/products/detail/phoenix-contact/3240084/3603842"),
  @("H18", "https://www.digikey.com/en/products/detail/cui-devices/DS04-254-2-03BK-SMT/11310920"),
  @("H15", "https://www.digikey.com/en/products/detail/molex/0022292051/1130591"),
  @("H16", "https://www.digikey.com/en/products/detail/molex/0010112053/171981"),
  @("H17", "https://www.digikey.com/en/products/detail/molex/0008550124/1784904"),
  @("H21", "https://www.digikey.com/en/products/detail/stackpole-electronics-inc/RMCP2010JT100R/2502808"),
  @("H9",  "https://www.digikey.com/en/products/detail/allegro-microsystems/A4956GESTR-T/5809983"),
  @("H22", "https://www.digikey.com/en/products/detail/microchip-technology/ATTINY4313-MMHR/3046321"),
  @("H19", "https://www.digikey.com/en/products/detail/vishay-dale-thin-film/THJP2512AST1/11313289?s=N4IgTCBcDaICoAkBSAFMBWAjGAggZTkxAF0BfIA"),
  @("H14", "https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/74650195R/6643984?s=N4IgTCBcDaIAQHYAsA2ArABgIwE40CUQBdAXyA"),
  @("H5",  "https://www.digikey.com/en/products/detail/avx-corporation/1812PC474KAT1A/1144366?s=N4IgTCBcDaIIwA45gAoGEAsB2DBpAggCpz4gC6AvkA"),
  @("H7",  "https://www.digikey.com/en/products/detail/vishay-dale/CRCW201015K0FKEF/1198590"),
  @("H4",  "https://www.digikey.com/en/products/detail/panasonic-electronic-components/EEE-HD1C472AM/9593462?s=N4IgTCBcDaIKIILQAkAiBGAwgFgOxgEEBZEAXQF8g"),
  @("H6",  "https://www.digikey.com/en/products/detail/avx-corporation/12101A621JAT2A/1605218?s=N4IgTCBcDaICwHYAcBaAjABgJwbelAcgCIgC6AvkA"),
  @("H10", "https://www.digikey.com/en/products/detail/maxim-integrated/MAX15062BATA-T/2591308"),
  @("H12", "https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/885012209005/5453563"),
  @("H13", "https://www.digikey.com/en/products/detail/holy-stone-enterprise-co-ltd/C1210X105K101T/13279995"),
  @("H3",  "https://www.digikey.com/en/products/detail/tdk-corporation/FG16X7R2A105KNT06/5811767"),
  @("H8",  "https://www.infineon.com/dgdl/irfs3306pbf.pdf?fileId=5546d462533600a40153563682652165")
)

foreach ($pair in $links) {
    $cellRef = $pair[0]
    $target = $pair[1]
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. View state: scrolled/selected cell moved from B15 to D9.
# ---------------------------------------------------------------------------
$ws.Range("D9").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
